$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '71.152.27'
$ws.Range('E2').Value = '  +3.53%  '

# Row 3
$ws.Range('D3').Value = '2.623.93'
$ws.Range('E3').Value = '  +3.51%  '

# Row 4
$ws.Range('E4').Value = '  +0.05%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '605.98'
$ws.Range('E5').Value = '  +1.94%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '180.95'
$ws.Range('E6').Value = '  +1.54%  '

# Row 7
$ws.Range('E7').Value = '  +0.00%  '

# Row 8
$ws.Range('E8').Value = '  +1.01%  '

# Row 9
$ws.Range('D9').Value = '2.622.57'
$ws.Range('E9').Value = '  +3.48%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.165'
$ws.Range('E10').Value = '  +13.30%  '

# Row 11
$ws.Range('E11').Value = '  +0.12%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.346'
$ws.Range('E12').Value = '  +2.14%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.01'
$ws.Range('E13').Value = '  +0.11%  '

# Row 14
$ws.Range('B14').Value = 'ShibaInu'
$ws.Range('C14').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000188'
$ws.Range('E14').Value = '  +9.53%  '

# Row 15
$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').Value = '3.074.89'
$ws.Range('E15').Value = '  +2.34%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '26.57'
$ws.Range('E16').Value = '  +1.78%  '

# Row 17
$ws.Range('D17').Value = '71.048.51'
$ws.Range('E17').Value = '  +3.76%  '

# Row 18
$ws.Range('D18').Value = '2.626.21'
$ws.Range('E18').Value = '  +4.80%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '383.48'
$ws.Range('E19').Value = '  +8.46%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.91'
$ws.Range('E20').Value = '  +5.18%  '

# Row 21
$ws.Range('E21').Value = '  +3.19%  '

# Row 22
$ws.Range('E22').Value = '  -1.80%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '72.11'
$ws.Range('E23').Value = '  +1.20%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.45'
$ws.Range('E24').Value = '  +5.42%  '

# Row 25
$ws.Range('E25').Value = '  -0.18%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.90'
$ws.Range('E26').Value = '  +10.12%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.64'
$ws.Range('E27').Value = '  +6.14%  '

# Row 28
$ws.Range('D28').Value = '2.758.26'
$ws.Range('E28').Value = '  +4.51%  '

# Row 29
$ws.Range('E29').Value = '  -0.27%  '

# Row 30
$ws.Range('D30').Value = '0.0₃0961'
$ws.Range('E30').Value = '  +6.61%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '540.13'
$ws.Range('E31').Value = '  +4.76%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.04'
$ws.Range('E32').Value = '  +2.71%  '

# Row 33
$ws.Range('E33').Value = '  +4.53%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.83'
$ws.Range('E34').Value = '  +3.05%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  -0.03%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '164.64'
$ws.Range('E36').Value = '  +0.24%  '

# Row 37
$ws.Range('E37').Value = '  -2.31%  '

# Row 38
$ws.Range('E38').Value = '  +4.07%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.90'
$ws.Range('E39').Value = '  +7.68%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '19.02'
$ws.Range('E40').Value = '  +1.70%  '

# Row 41
$ws.Range('E41').Value = '  +4.58%  '

# Row 42
$ws.Range('E42').Value = '  +8.60%  '

# Row 43
$ws.Range('E43').Value = '  +0.03%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.04'
$ws.Range('E44').Value = '  +3.82%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.330'
$ws.Range('E45').Value = '  +0.85%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '39.91'
$ws.Range('E46').Value = '  +2.32%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '153.89'
$ws.Range('E47').Value = '  +0.63%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.63'
$ws.Range('E48').Value = '  +1.55%  '

# Row 49
$ws.Range('E49').Value = '  +4.66%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.532'
$ws.Range('E50').Value = '  +2.07%  '

# Row 51
$ws.Range('D51').Value = '0.0₆0261'
$ws.Range('E51').Value = '  +0.01%  '
